$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds text-formatted numeric-looking strings
# (e.g. thousand-separated prices, trailing zeros). Force text format
# on the whole column before assigning so Excel does not silently
# reinterpret/normalize the values as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.465.70"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.569.66"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "290.53"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "0.3692"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").Value = "0.3372"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").Value = "1.147"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("D11").Value = "0.07526"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "6.016"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").Value = "6.962"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "1.570.20"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "0.00001120"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "90.43"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").Value = "0.06773"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "6.350"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").Value = "16.39"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").Value = "12.23"
$ws.Range("E23").Value = "  +2.77%  "
$ws.Range("D24").Value = "22.456.17"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "2.381"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "2.648"
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("D27").Value = "20.02"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "148.96"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("D29").Value = "5.054"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").Value = "124.96"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "1.750.56"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").Value = "1.064"
$ws.Range("E32").Value = "  +8.19%  "
$ws.Range("D33").Value = "6.190"
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("D34").Value = "2.012"
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("D35").Value = "9.801"
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("D36").Value = "0.08349"
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("D37").Value = "0.02471"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("D38").Value = "1.352"
$ws.Range("E38").Value = "  -5.11%  "
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").Value = "0.06558"
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("D41").Value = "5.428"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").Value = "11.21"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").Value = "0.6214"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").Value = "14.15"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "3.807"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "0.5851"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D50").Value = "1.238"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").Value = "0.07305"
$ws.Range("E51").Value = "  -0.16%  "

# Rows 48/49: NEARProtocol and Quant swapped positions (with updated figures)
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "2.067"
$ws.Range("E48").Value = "  +0.61%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "127.96"
$ws.Range("E49").Value = "  +2.81%  "
